# EindeStemperiode.xlsx update
# - Adds a "Bron" (source) link for the 2009 row (C13).
# - Corrects the end-of-voting-period dates for 2005, 2004 and 2003 (A17:A19)
#   and adds sources for those years (C17:C19).
# - Adds "Gok" (guess) as the source for 2002, 2001, 2000 and 1999 (C20:C23).
# - Moves the active selection to C4, matching the author's final cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# New source for 2009 (row 13)
$ws.Range("C13").Value = "https://www.parool.nl/nieuws/veel-animo-voor-stemmen-op-top-2000~b743b9a8/"

# Corrected dates for 2005 / 2004 / 2003
$ws.Range("A17").Value = Get-Date -Year 2005 -Month 11 -Day 30 -Hour 0 -Minute 0 -Second 0
$ws.Range("A18").Value = Get-Date -Year 2004 -Month 12 -Day 10 -Hour 0 -Minute 0 -Second 0
$ws.Range("A19").Value = Get-Date -Year 2003 -Month 12 -Day 15 -Hour 0 -Minute 0 -Second 0

# Sources for 2004, 2005, 2003 -- set in this order so new shared strings
# are appended in the same sequence as in the reference workbook
$ws.Range("C18").Value = "https://www.radiofreak.nl/stemmen-voor-top-2000-begint-morgen/"
$ws.Range("C17").Value = "https://www.radiofreak.nl/stemming-voor-radio-2-top-2000-begonnen/"
$ws.Range("C19").Value = "https://www.radiofreak.nl/radio-2-start-maandag-met-de-top-2000/"

# 2002, 2001, 2000, 1999 (rows 20-23): mark source as a guess ("Gok")
$ws.Range("C20").Value = "Gok"
$ws.Range("C21").Value = "Gok"
$ws.Range("C22").Value = "Gok"
$ws.Range("C23").Value = "Gok"

# Restore the author's active cell selection
$ws.Range("C4").Select()
